$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style of the other headers (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"
$excel.CutCopyMode = 0

# Populate the time_taken values for each data row (F2:F15)
$timestamps = @(
    "2021-10-05 13:40:37.664731",
    "2021-10-05 13:40:37.664742",
    "2021-10-05 13:40:37.664745",
    "2021-10-05 13:40:37.664748",
    "2021-10-05 13:40:37.664750",
    "2021-10-05 13:40:37.664753",
    "2021-10-05 13:40:37.664755",
    "2021-10-05 13:40:37.664758",
    "2021-10-05 13:40:37.664760",
    "2021-10-05 13:40:37.664763",
    "2021-10-05 13:40:37.664765",
    "2021-10-05 13:40:37.664768",
    "2021-10-05 13:40:37.664770",
    "2021-10-05 13:40:37.664773"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
